# Adds columns_experiment rows 9 and 10 (spreadsheet rows 10 and 11) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 10 (experiment #9) ----
$ws.Range("B10").Value = "VGG19`n+ Dense(512, relu, regularizer)`n+ Dropout`n+ Dense(512, relu, regularizer)`n+ Dropout"
$ws.Range("C10").Value = "Shift [-4, 4] (mode='edge')`nrotate [-0.5, 0.5]"
$ws.Range("D10").Value = "Adjust param according to the random crop`nGrayscale`nUse border between windows`nCrop sky and shop"
$ws.Range("E10").Value = "x20"
$ws.Range("F10").Value = 0.0001
$ws.Range("G10").Value = 60
$ws.Range("H10").Value = 0.0101
$ws.Range("I10").Value = 0.00017908
$ws.Range("J10").Value = 0.0234
$ws.Range("K10").Value = 0.0035

# ---- Row 11 (experiment #10) ----
$ws.Range("B11").Value = "VGG19`n+ Dense(512, relu, regularizer)`n+ Dropout`n+ Dense(512, relu, regularizer)`n+ Dropout"
$ws.Range("C11").Value = "Shift [-4, 4] (mode='edge')`nrotate [-0.5, 0.5] (mode='nearest')"
$ws.Range("D11").Value = "Adjust param according to the random crop`nGrayscale`nUse border between windows`nCrop sky and shop"
$ws.Range("E11").Value = "x20"
$ws.Range("F11").Value = 0.0001
$ws.Range("G11").Value = 60
$ws.Range("H11").Value = 0.0097
$ws.Range("I11").Value = 0.00016637
$ws.Range("J11").Value = 0.0306
$ws.Range("K11").Value = 0.007

"done"
